$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "reena"
$ws.Range("A2").Value = "nidhi"
$ws.Range("A3").Value = "mokshi"
$ws.Range("A4").Value = "nikita"

$ws.Range("A5").Select()
